$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data row for the SNSWLinkingTestB row (row 4) - set first so the shared
# string for the address is created before the new header strings.
$ws.Range("M4").Value = "7 KEVIN RD, ALBION PARK NSW 2527"
$ws.Range("N4").Value = "7 KEVIN RD, ALBION PARK NSW 2527"

# New header columns for Option B Linking Details
$ws.Range("N1").Value = "Mailing Address"
$ws.Range("M1").Value = "Residential Adress"

# Match header style (yellow fill) used by the rest of row 1
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Best-fit the new columns to their (long) contents, like Excel's
# "AutoFit Column Width" on columns M:N.
$ws.Range("M1:N1").ColumnWidth = 32.59

# Update selection to match target state
$ws.Range("K10:K11").Select()
